$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2759.6155
$ws.Range("I17").Value = 3411
$ws.Range("J17").Value = 1999.6666
$ws.Range("K17").Value = 10233
$ws.Range("L17").Value = 5998.9998
$ws.Range("M17").Value = -10065
$ws.Range("N17").Value = -6334.9998

$ws.Range("H47").Value = 13859.8
$ws.Range("I47").Value = 13859.8
$ws.Range("K47").Value = 13859.8
$ws.Range("M47").Value = -12887.8

$ws.Range("H54").Value = 9970
$ws.Range("I54").Value = 9970
$ws.Range("K54").Value = 9970
$ws.Range("M54").Value = -9484

$ws.Range("H118").Value = 209.5
$ws.Range("I118").Value = 209.5
$ws.Range("K118").Value = 628.5
$ws.Range("M118").Value = 1028.5

$ws.Range("H131").Value = 4459
$ws.Range("I131").Value = 1361.125
$ws.Range("J131").Value = 7999.4287
$ws.Range("K131").Value = 4083.375
$ws.Range("L131").Value = 23998.2861
$ws.Range("M131").Value = 956.625
$ws.Range("N131").Value = -34078.2861

$ws.Range("H137").Value = 2439.9412
$ws.Range("I137").Value = 964.7778
$ws.Range("K137").Value = 2894.3334
$ws.Range("M137").Value = -344.3334

$ws.Range("H138").Value = 2848.186
$ws.Range("I138").Value = 1666.6666
$ws.Range("J138").Value = 2936.8
$ws.Range("K138").Value = 4999.9998
$ws.Range("L138").Value = 8810.400000000001
$ws.Range("M138").Value = 140.0002000000004
$ws.Range("N138").Value = -19090.4

$ws.Range("H141").Value = 2304.3635
$ws.Range("I141").Value = 1534.8
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 4604.4
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = 575.6000000000004
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2186.375
$ws.Range("I45").Value = 2198.5
$ws.Range("J45").Value = 2150
$ws.Range("K45").Value = 2198.5
$ws.Range("L45").Value = 2150
$ws.Range("M45").Value = -1821.5
$ws.Range("N45").Value = -2904

$ws.Range("H61").Value = 1942.125
$ws.Range("I61").Value = 1942.125
$ws.Range("K61").Value = 1942.125
$ws.Range("M61").Value = -1730.125

$ws.Range("H88").Value = 380.66666
$ws.Range("I88").Value = 423.44446
$ws.Range("K88").Value = 423.44446
$ws.Range("M88").Value = -17.44445999999999

$ws.Range("H91").Value = 380.66666
$ws.Range("I91").Value = 423.44446
$ws.Range("K91").Value = 423.44446
$ws.Range("M91").Value = 980.5555400000001

$ws.Range("H95").Value = 42997.8
$ws.Range("J95").Value = 42997.8
$ws.Range("L95").Value = 42997.8
$ws.Range("N95").Value = -48489.8

$ws.Range("H122").Value = 558786.1
$ws.Range("I122").Value = 1112687.4
$ws.Range("J122").Value = 4884.8887
$ws.Range("K122").Value = 3338062.2
$ws.Range("L122").Value = 14654.6661
$ws.Range("M122").Value = -3335612.2
$ws.Range("N122").Value = -19554.6661

$ws.Range("H132").Value = 3971.1
$ws.Range("I132").Value = 2963.875
$ws.Range("K132").Value = 8891.625
$ws.Range("M132").Value = -6361.625

$ws.Range("H136").Value = 1942.125
$ws.Range("I136").Value = 1942.125
$ws.Range("K136").Value = 5826.375
$ws.Range("M136").Value = -3276.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2453.2593
$ws.Range("I99").Value = 2197.625
$ws.Range("K99").Value = 2197.625
$ws.Range("M99").Value = -699.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 25200.334
$ws.Range("J92").Value = 25200.334
$ws.Range("L92").Value = 25200.334
$ws.Range("N92").Value = -30192.334

$ws.Range("H132").Value = 2256.3
$ws.Range("I132").Value = 1799.375
$ws.Range("K132").Value = 5398.125
$ws.Range("M132").Value = -2868.125

$ws.Range("H141").Value = 20508.945
$ws.Range("J141").Value = 20508.945
$ws.Range("L141").Value = 20508.945
$ws.Range("N141").Value = -30868.945

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 66750668
$ws.Range("I4").Value = 87065416
$ws.Range("K4").Value = 261196248
$ws.Range("M4").Value = -261196136

$ws.Range("H14").Value = 874.4545000000001
$ws.Range("I14").Value = 874.4545000000001
$ws.Range("K14").Value = 2623.3635
$ws.Range("M14").Value = -2450.3635

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H29").Value = 1166
$ws.Range("I29").Value = 1000
$ws.Range("K29").Value = 3000
$ws.Range("M29").Value = -2723

$ws.Range("H107").Value = 669.7714
$ws.Range("J107").Value = 646.21875
$ws.Range("L107").Value = 1938.65625
$ws.Range("N107").Value = -5778.65625

$ws.Range("H113").Value = 1135.1578
$ws.Range("J113").Value = 1051.7693
$ws.Range("L113").Value = 3155.3079
$ws.Range("N113").Value = -7495.3079

$ws.Range("H139").Value = 4798.1113
$ws.Range("I139").Value = 2531.3333
$ws.Range("K139").Value = 7593.999899999999
$ws.Range("M139").Value = -2453.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 900
$ws.Range("I80").Value = 900
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 900
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 98
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 900
$ws.Range("I83").Value = 900
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 492
$ws.Range("N83").ClearContents()

$ws.Range("H102").Value = 1586.875
$ws.Range("I102").Value = 426.13043
$ws.Range("K102").Value = 426.13043
$ws.Range("M102").Value = 1195.86957

$ws.Range("H107").Value = 855.86664
$ws.Range("I107").Value = 902.0625
$ws.Range("J107").Value = 803.0714
$ws.Range("K107").Value = 902.0625
$ws.Range("L107").Value = 803.0714
$ws.Range("M107").Value = 1017.9375
$ws.Range("N107").Value = -4643.0714

$ws.Range("H126").Value = 3677.25
$ws.Range("I126").Value = 1906
$ws.Range("J126").Value = 4942.4287
$ws.Range("K126").Value = 5718
$ws.Range("L126").Value = 14827.2861
$ws.Range("M126").Value = -3248
$ws.Range("N126").Value = -19767.2861

$ws.Range("H132").Value = 2719
$ws.Range("I132").Value = 1861.2142
$ws.Range("K132").Value = 5583.642599999999
$ws.Range("M132").Value = -3053.642599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2323.25
$ws.Range("I7").Value = 2226.7144
$ws.Range("K7").Value = 2226.7144
$ws.Range("M7").Value = -2114.7144

$ws.Range("H46").Value = 3308.7083
$ws.Range("I46").Value = 2442.111
$ws.Range("K46").Value = 2442.111
$ws.Range("M46").Value = -2254.111

$ws.Range("H122").Value = 12998.2
$ws.Range("I122").Value = 13748.5
$ws.Range("J122").Value = 9997
$ws.Range("K122").Value = 41245.5
$ws.Range("L122").Value = 29991
$ws.Range("M122").Value = -38795.5
$ws.Range("N122").Value = -34891

$ws.Range("H126").Value = 2323.25
$ws.Range("I126").Value = 2226.7144
$ws.Range("K126").Value = 6680.1432
$ws.Range("M126").Value = -4210.1432

$ws.Range("H132").Value = 3373.682
$ws.Range("I132").Value = 2797.0278
$ws.Range("K132").Value = 8391.0834
$ws.Range("M132").Value = -5861.0834

$ws.Range("H136").Value = 1637.7858
$ws.Range("I136").Value = 1630
$ws.Range("K136").Value = 4890
$ws.Range("M136").Value = -2340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 732.5
$ws.Range("I107").Value = 404.16666
$ws.Range("K107").Value = 1212.49998
$ws.Range("M107").Value = 707.5000199999999

$ws.Range("H122").Value = 16990
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 16990
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 50970
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -55870

$ws.Range("H126").Value = 2939
$ws.Range("I126").Value = 864.6
$ws.Range("J126").Value = 8125
$ws.Range("K126").Value = 2593.8
$ws.Range("L126").Value = 24375
$ws.Range("M126").Value = -123.8000000000002
$ws.Range("N126").Value = -29315

$ws.Range("H132").Value = 1955.7
$ws.Range("I132").Value = 1383.375
$ws.Range("K132").Value = 4150.125
$ws.Range("M132").Value = -1620.125

$ws.Range("H136").Value = 4407.3335
$ws.Range("I136").Value = 1199.125
$ws.Range("K136").Value = 3597.375
$ws.Range("M136").Value = -1047.375
